# Auto-generated edit script: updates cryptocurrency price/volume/name/link
# cells in Sheet1 to match the "Updated symbol list" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue([string]$cellRef, [string]$text) {
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextValue "D2" "243.21"
Set-TextValue "D4" "5.402"
Set-TextValue "D5" "0.05995"
Set-TextValue "D6" "3.427"
Set-TextValue "D7" "6.514"
Set-TextValue "D8" "0.8116"
Set-TextValue "D9" "0.9215"
Set-TextValue "B10" "One"
Set-TextValue "C10" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D10" "0.01120"
Set-TextValue "E10" "9OneONEBestin24h"
Set-TextValue "B11" "WazirX"
Set-TextValue "C11" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D11" "0.1428"
Set-TextValue "E11" "10WazirXWRX"
Set-TextValue "B12" "MandalaExchangeToken"
Set-TextValue "C12" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D12" "0.07424"
Set-TextValue "E12" "11MandalaExchangeTokenMDX"
Set-TextValue "B13" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C13" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D13" "0.03304"
Set-TextValue "E13" "12LiechtensteinCryptoassetsExchangeLCX"
Set-TextValue "B14" "BitrueCoin"
Set-TextValue "C14" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D14" "0.03066"
Set-TextValue "E14" "13BitrueCoinBTR"
Set-TextValue "B15" "BitMartToken"
Set-TextValue "C15" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D15" "0.09353"
Set-TextValue "E15" "14BitMartTokenBMX"
Set-TextValue "B16" "MCDex"
Set-TextValue "C16" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D16" "3.851"
Set-TextValue "E16" "15MCDexMCB"
Set-TextValue "B17" "BitForexToken"
Set-TextValue "C17" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D17" "0.001578"
Set-TextValue "E17" "16BitForexTokenBF"
Set-TextValue "B18" "CoinExToken"
Set-TextValue "C18" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D18" "0.04696"
Set-TextValue "E18" "17CoinExTokenCET"
Set-TextValue "D19" "0.005854"
Set-TextValue "D21" "0.004880"
Set-TextValue "D23" "3.567"
Set-TextValue "D24" "2.137"
Set-TextValue "D25" "0.3234"
Set-TextValue "D26" "0.1331"
Set-TextValue "D40" "0.03968"
Set-TextValue "B41" "KickToken"
Set-TextValue "C41" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.006381"
Set-TextValue "E41" "40KickTokenKICK"
Set-TextValue "B42" "CEJI"
Set-TextValue "C42" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.003800"
Set-TextValue "E42" "41CEJICEJI"
Set-TextValue "D44" "0.009186"
Set-TextValue "D45" "0.00005076"
Set-TextValue "D47" "0.7001"
Set-TextValue "E47" "46CoinbaseStockTokenCOINWorstin24h"
